# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
#
# Every existing data cell in this workbook (dates, clock times, hour
# buckets, percentages, temperatures, ...) is stored as literal TEXT, never
# as a native Excel number/date. Plain `Range.Value = "..."` lets Excel's
# COM layer auto-detect dates ("2026-02-01") and percentages ("78.3%") and
# silently coerce them into numeric/date serials, which would diverge from
# the source file's inline-string cells. To avoid that we:
#   1. Force the target block's NumberFormat to Text ("@") before writing,
#      so the assignment is stored verbatim as a string.
#   2. After writing, reset the block's Style back to "Normal" so no stray
#      number-format style lingers on the cells (matching the original
#      file, where none of these cells carry an explicit style index).

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param($Worksheet, $Rows)

    $firstRow = $Rows[0][0]
    $lastRow = $Rows[$Rows.Count - 1][0]

    $targetRange = $Worksheet.Range("A$($firstRow):F$($lastRow)")
    $targetRange.NumberFormat = "@"

    foreach ($row in $Rows) {
        $r = $row[0]
        $Worksheet.Cells.Item($r, 1).Value = $row[1]
        $Worksheet.Cells.Item($r, 2).Value = $row[2]
        $Worksheet.Cells.Item($r, 3).Value = $row[3]
        $Worksheet.Cells.Item($r, 4).Value = $row[4]
        $Worksheet.Cells.Item($r, 5).Value = $row[5]
        $Worksheet.Cells.Item($r, 6).Value = $row[6]
    }

    $targetRange.Style = "Normal"
}

# --- PIR: append rows 69-82 -------------------------------
$pirRows = @(
    @(69,'2026-02-01','19:58:11','19:00','Bathroom','No Motion','Inactive'),
    @(70,'2026-02-01','19:58:12','19:00','Bathroom','No Motion','Inactive'),
    @(71,'2026-02-01','19:58:15','19:00','Bathroom','No Motion','Inactive'),
    @(72,'2026-02-01','19:58:20','19:00','Bathroom','No Motion','Inactive'),
    @(73,'2026-02-01','19:58:25','19:00','Bathroom','No Motion','Inactive'),
    @(74,'2026-02-01','19:58:30','19:00','Bathroom','No Motion','Inactive'),
    @(75,'2026-02-01','19:58:35','19:00','Bathroom','No Motion','Inactive'),
    @(76,'2026-02-01','19:58:40','19:00','Bathroom','No Motion','Inactive'),
    @(77,'2026-02-01','19:58:45','19:00','Bathroom','No Motion','Inactive'),
    @(78,'2026-02-01','19:58:50','19:00','Bathroom','No Motion','Inactive'),
    @(79,'2026-02-01','19:58:51','19:00','Bathroom','Motion Detected','Active'),
    @(80,'2026-02-01','19:58:58','19:00','Bathroom','No Motion','Inactive'),
    @(81,'2026-02-01','19:59:03','19:00','Bathroom','No Motion','Inactive'),
    @(82,'2026-02-01','19:59:08','19:00','Bathroom','No Motion','Inactive')
)
$pirSheet = $wb.Worksheets.Item("PIR")
Append-Rows $pirSheet $pirRows

# --- Humidity: append rows 56-64 --------------------------
$humidityRows = @(
    @(56,'2026-02-01','19:58:11','19:00','Bathroom','78.3%','Active'),
    @(57,'2026-02-01','19:58:12','19:00','Bathroom','76.9%','Active'),
    @(58,'2026-02-01','19:58:17','19:00','Bathroom','78.1%','Active'),
    @(59,'2026-02-01','19:58:22','19:00','Bathroom','77.0%','Active'),
    @(60,'2026-02-01','19:58:32','19:00','Bathroom','77.3%','Active'),
    @(61,'2026-02-01','19:58:37','19:00','Bathroom','78.4%','Active'),
    @(62,'2026-02-01','19:58:57','19:00','Bathroom','77.8%','Active'),
    @(63,'2026-02-01','19:59:02','19:00','Bathroom','77.4%','Active'),
    @(64,'2026-02-01','19:59:07','19:00','Bathroom','78.5%','Active')
)
$humiditySheet = $wb.Worksheets.Item("Humidity")
Append-Rows $humiditySheet $humidityRows

# --- Temperature: append rows 56-64 ---------------------
$temperatureRows = @(
    @(56,'2026-02-01','19:58:11','19:00','Bathroom','25.2C','Active'),
    @(57,'2026-02-01','19:58:12','19:00','Bathroom','25.1C','Active'),
    @(58,'2026-02-01','19:58:17','19:00','Bathroom','25.1C','Active'),
    @(59,'2026-02-01','19:58:22','19:00','Bathroom','25.2C','Active'),
    @(60,'2026-02-01','19:58:32','19:00','Bathroom','25.1C','Active'),
    @(61,'2026-02-01','19:58:37','19:00','Bathroom','25.1C','Active'),
    @(62,'2026-02-01','19:58:57','19:00','Bathroom','25.1C','Active'),
    @(63,'2026-02-01','19:59:02','19:00','Bathroom','25.1C','Active'),
    @(64,'2026-02-01','19:59:07','19:00','Bathroom','25.1C','Active')
)
$temperatureSheet = $wb.Worksheets.Item("Temperature")
Append-Rows $temperatureSheet $temperatureRows
